$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextValue $ws "D2" "305.72"
Set-TextValue $ws "E2" "-0.52%"
Set-TextValue $ws "G2" "22"

Set-TextValue $ws "D3" "38.79"
Set-TextValue $ws "E3" "6.39%"
Set-TextValue $ws "G3" "22"

Set-TextValue $ws "D4" "5.115"
Set-TextValue $ws "E4" "0.86%"
Set-TextValue $ws "G4" "22"

Set-TextValue $ws "D5" "0.08074"
Set-TextValue $ws "E5" "-0.21%"
Set-TextValue $ws "G5" "22"

Set-TextValue $ws "E6" "-3.91%"
Set-TextValue $ws "G6" "22"

Set-TextValue $ws "E7" "0.63%"
Set-TextValue $ws "G7" "22"

Set-TextValue $ws "D8" "8.012"
Set-TextValue $ws "G8" "22"

Set-TextValue $ws "D9" "0.9274"
Set-TextValue $ws "E9" "-0.08%"
Set-TextValue $ws "G9" "22"

Set-TextValue $ws "D10" "0.1434"
Set-TextValue $ws "E10" "-2.31%"
Set-TextValue $ws "G10" "22"

Set-TextValue $ws "D11" "0.1914"
Set-TextValue $ws "E11" "-1.56%"
Set-TextValue $ws "G11" "22"

Set-TextValue $ws "D12" "0.09047"
Set-TextValue $ws "E12" "-0.86%"
Set-TextValue $ws "G12" "22"

Set-TextValue $ws "D13" "0.03511"
Set-TextValue $ws "E13" "-0.22%"
Set-TextValue $ws "G13" "22"

Set-TextValue $ws "D14" "0.09777"
Set-TextValue $ws "E14" "-1.12%"
Set-TextValue $ws "G14" "22"

Set-TextValue $ws "D15" "0.001394"
Set-TextValue $ws "E15" "-0.94%"
Set-TextValue $ws "G15" "22"

Set-TextValue $ws "D16" "0.005919"
Set-TextValue $ws "E16" "-7.35%"
Set-TextValue $ws "G16" "22"

Set-TextValue $ws "D17" "3.775"
Set-TextValue $ws "E17" "-1.70%"
Set-TextValue $ws "G17" "22"

Set-TextValue $ws "E18" "-3.14%"
Set-TextValue $ws "G18" "22"

Set-TextValue $ws "E19" "0.29%"
Set-TextValue $ws "G19" "22"

Set-TextValue $ws "D20" "0.1327"
Set-TextValue $ws "E20" "-0.06%"
Set-TextValue $ws "G20" "22"

Set-TextValue $ws "D21" "4.685"
Set-TextValue $ws "E21" "-2.85%"
Set-TextValue $ws "G21" "22"

Set-TextValue $ws "E22" "2.96%"
Set-TextValue $ws "G22" "22"

Set-TextValue $ws "D23" "0.04373"
Set-TextValue $ws "E23" "-0.38%"
Set-TextValue $ws "G23" "22"

Set-TextValue $ws "D24" "0.001226"
Set-TextValue $ws "E24" "-0.92%"
Set-TextValue $ws "G24" "22"

Set-TextValue $ws "D25" "0.004268"
Set-TextValue $ws "E25" "2.08%"
Set-TextValue $ws "G25" "22"

Set-TextValue $ws "E26" "-0.04%"
Set-TextValue $ws "G26" "22"

Set-TextValue $ws "G27" "22"

Set-TextValue $ws "G28" "22"

Set-TextValue $ws "G29" "22"

Set-TextValue $ws "G30" "22"

Set-TextValue $ws "G31" "22"

Set-TextValue $ws "G32" "22"

Set-TextValue $ws "G33" "22"

Set-TextValue $ws "G34" "22"

Set-TextValue $ws "G35" "22"

Set-TextValue $ws "G36" "22"

Set-TextValue $ws "G37" "22"

Set-TextValue $ws "G38" "22"

Set-TextValue $ws "D39" "0.02026"
Set-TextValue $ws "E39" "-0.97%"
Set-TextValue $ws "G39" "22"

Set-TextValue $ws "D40" "0.05044"
Set-TextValue $ws "E40" "-1.70%"
Set-TextValue $ws "G40" "22"

Set-TextValue $ws "D41" "0.007523"
Set-TextValue $ws "E41" "0.59%"
Set-TextValue $ws "G41" "22"

Set-TextValue $ws "D42" "0.009776"
Set-TextValue $ws "E42" "-3.10%"
Set-TextValue $ws "G42" "22"

Set-TextValue $ws "D43" "0.1340"
Set-TextValue $ws "E43" "-2.18%"
Set-TextValue $ws "G43" "22"

Set-TextValue $ws "D44" "0.002105"
Set-TextValue $ws "E44" "-0.98%"
Set-TextValue $ws "G44" "22"

Set-TextValue $ws "D45" "0.009900"
Set-TextValue $ws "E45" "0.02%"
Set-TextValue $ws "G45" "22"

Set-TextValue $ws "D46" "0.00006218"
Set-TextValue $ws "E46" "-1.56%"
Set-TextValue $ws "G46" "22"

Set-TextValue $ws "E47" "-0.08%"
Set-TextValue $ws "G47" "22"

Set-TextValue $ws "D48" "0.002876"
Set-TextValue $ws "G48" "22"

Set-TextValue $ws "E49" "12.41%"
Set-TextValue $ws "G49" "22"

Set-TextValue $ws "E50" "-0.08%"
Set-TextValue $ws "G50" "22"

Set-TextValue $ws "E51" "-0.08%"
Set-TextValue $ws "G51" "22"

